$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the explanatory note in G16 (new shared string) ---
$ws.Range("G16").Value = "These points is removed due to Martin Jørgensen's discovery of the error"

# --- Move the old row17/row18 data (A:E) to G:K on rows 17/18 ---
# Old row17: 2020.58, -50, 50, 5, 5
$ws.Range("G17").Value = 2020.58
$ws.Range("H17").Value = -50
$ws.Range("I17").Value = 50
$ws.Range("J17").Value = 5
$ws.Range("K17").Value = 5

# Old row18: 2021.24, -51, 34, 0.5, 0.5
$ws.Range("G18").Value = 2021.24
$ws.Range("H18").Value = -51
$ws.Range("I18").Value = 34
$ws.Range("J18").Value = 0.5
$ws.Range("K18").Value = 0.5

# --- Shift A:E data of rows 19,20,21 up into rows 17,18,19 ---
# Old row19 -> new row17 (A:E)
$ws.Range("A17").Value = 2021.41
$ws.Range("B17").Value = 2.9
$ws.Range("C17").Value = -11.76
$ws.Range("D17").Value = 0.5
$ws.Range("E17").Value = 0.5

# Old row20 -> new row18 (A:E)
$ws.Range("A18").Value = 2021.47
$ws.Range("B18").Value = 12.9
$ws.Range("C18").Value = -2.9
$ws.Range("D18").Value = 0.5
$ws.Range("E18").Value = 0.5

# Old row21 -> new row19 (A:E)
$ws.Range("A19").Value = 2021.56
$ws.Range("B19").Value = 18.2
$ws.Range("C19").Value = 9.41
$ws.Range("D19").Value = 0.5
$ws.Range("E19").Value = 0.5

# --- Clear now-unused rows 20:21 (A:E) ---
$ws.Range("A20:E21").ClearContents()

# --- Update selection to match target view ---
$ws.Range("E23").Select() | Out-Null

# --- Page setup (paperSize/orientation) ---
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

Write-Host "done"
